$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number need to be forced to
# remain text (matching the original inline-string cell type) instead of being
# auto-converted to a numeric cell by Excel.
$textCells = @('D5', 'D6', 'D10', 'D11', 'D12', 'D13', 'D14', 'D20', 'D21', 'D23', 'D24', 'D25', 'D27', 'D28', 'D31', 'D32', 'D33', 'D34', 'D36', 'D39', 'D40', 'D41', 'D43', 'D44', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '43.854.96'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').Value = '2.295.37'
$ws.Range('E3').Value = '  +0.03%  '
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('D5').Value = '113.15'
$ws.Range('E5').Value = '  +16.94%  '
$ws.Range('D6').Value = '269.86'
$ws.Range('E6').Value = '  +0.48%  '
$ws.Range('E7').Value = '  +0.49%  '
$ws.Range('E8').Value = '  +0.29%  '
$ws.Range('E9').Value = '  +2.28%  '
$ws.Range('D10').Value = '48.11'
$ws.Range('E10').Value = '  +6.19%  '
$ws.Range('D11').Value = '0.0952'
$ws.Range('E11').Value = '  +2.01%  '
$ws.Range('D12').Value = '9.27'
$ws.Range('E12').Value = '  +17.34%  '
$ws.Range('D13').Value = '0.107'
$ws.Range('E13').Value = '  +0.38%  '
$ws.Range('D14').Value = '15.82'
$ws.Range('E14').Value = '  +0.88%  '
$ws.Range('D15').Value = '2.637.93'
$ws.Range('E15').Value = '  -0.11%  '
$ws.Range('E16').Value = '  -0.30%  '
$ws.Range('D17').Value = '2.297.58'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('D18').Value = '43.733.21'
$ws.Range('E18').Value = '  +0.06%  '
$ws.Range('E19').Value = '  -0.99%  '
$ws.Range('D20').Value = '6.75'
$ws.Range('E20').Value = '  +8.90%  '
$ws.Range('D21').Value = '72.34'
$ws.Range('E21').Value = '  +0.33%  '
$ws.Range('E22').Value = '  -2.29%  '
$ws.Range('B23').Value = 'InternetComputer(DFINITY)'
$ws.Range('C23').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D23').Value = '9.89'
$ws.Range('E23').Value = '  +8.41%  '
$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').Value = '232.69'
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('D25').Value = '2.84'
$ws.Range('E25').Value = '  +7.76%  '
$ws.Range('D27').Value = '11.71'
$ws.Range('E27').Value = '  +3.99%  '
$ws.Range('D28').Value = '41.94'
$ws.Range('E28').Value = '  +8.51%  '
$ws.Range('E29').Value = '  -2.08%  '
$ws.Range('E30').Value = '  -2.54%  '
$ws.Range('D31').Value = '175.39'
$ws.Range('E31').Value = '  +0.14%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '0.0929'
$ws.Range('E32').Value = '  +2.77%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '21.56'
$ws.Range('E33').Value = '  -1.16%  '
$ws.Range('D34').Value = '5.68'
$ws.Range('E34').Value = '  +4.78%  '
$ws.Range('E35').Value = '  +1.80%  '
$ws.Range('D36').Value = '4.68'
$ws.Range('E36').Value = '  +2.69%  '
$ws.Range('E37').Value = '  +3.26%  '
$ws.Range('E38').Value = '  +0.81%  '
$ws.Range('D39').Value = '3.83'
$ws.Range('E39').Value = '  +13.41%  '
$ws.Range('D40').Value = '74.61'
$ws.Range('E40').Value = '  +15.44%  '
$ws.Range('D41').Value = '13.81'
$ws.Range('E41').Value = '  +13.28%  '
$ws.Range('E42').Value = '  +2.43%  '
$ws.Range('D43').Value = '2.39'
$ws.Range('E43').Value = '  +2.92%  '
$ws.Range('D44').Value = '6.34'
$ws.Range('E44').Value = '  +22.90%  '
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('E46').Value = '  +3.71%  '
$ws.Range('D47').Value = '8.80'
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('D48').Value = '103.16'
$ws.Range('E48').Value = '  +6.14%  '
$ws.Range('D49').Value = '0.0996'
$ws.Range('E49').Value = '  -2.36%  '
$ws.Range('D50').Value = '1.23'
$ws.Range('E50').Value = '  +2.91%  '
$ws.Range('D51').Value = '0.465'
$ws.Range('E51').Value = '  +5.37%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
